# Code studies oct 2020
# Update BASE AMOUNT (F), INITIAL AMOUNT (G) and TOTAL (H) figures for rows 2-22
# on the "Almal" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Almal")

$values = @{
    2  = @(600, 360, 960)
    3  = @(1000, 0, 2000)
    4  = @(1000, 360, 2360)
    5  = @(1000, 360, 2360)
    6  = @(600, 360, 960)
    7  = @(600, 360, 960)
    8  = @(1000, 0, 2000)
    9  = @(600, 0, 600)
    10 = @(1000, 0, 1800)
    11 = @(600, 0, 600)
    12 = @(1000, 0, 2000)
    13 = @(600, 0, 600)
    14 = @(1000, 0, 2000)
    15 = @(600, 0, 600)
    16 = @(1000, 0, 2000)
    17 = @(600, 0, 600)
    18 = @(1000, 0, 2000)
    19 = @(600, 0, 600)
    20 = @(1000, 0, 2000)
    21 = @(600, 0, 600)
    22 = @(1000, 0, 2000)
}

foreach ($row in $values.Keys) {
    $triple = $values[$row]
    $ws.Cells.Item($row, 6).Value = $triple[0]
    $ws.Cells.Item($row, 7).Value = $triple[1]
    $ws.Cells.Item($row, 8).Value = $triple[2]
}
